# fix(publipostage): Try to solve Excel emoji problem
#
# The "statut" column (column A) used emoji glyphs (📕/📘/📙/📗) as status
# codes. Replace them with plain-text equivalents that render reliably:
#   📕 -> -3
#   📘 -> ⚠️
#   📙 -> +3
#   📗 -> ✅

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "📕" = "-3"
    "📘" = "⚠️"
    "📙" = "+3"
    "📗" = "✅"
}

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    if ($map.ContainsKey($old)) {
        $new = $map[$old]
        # "-3" and "+3" look like numbers to Excel, which would silently
        # convert them to numeric values (and drop the "+"). Prefix with an
        # apostrophe (exactly like typing '-3 into a cell) so they are
        # stored as text, matching the other two (already non-numeric)
        # replacement labels.
        if ($new -eq "-3" -or $new -eq "+3") {
            $cell.Value2 = "'" + $new
        } else {
            $cell.Value2 = $new
        }
    }
}
